$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header table (rows 2-4): gear-motor options ---
$ws.Range("B3").Value = 60

# --- "Best" motor label (row 7, merged A7:D7) and removal of the old
#     "9.8 meters" annotation on row 35 (E35). Clearing E35 first frees
#     the shared string slot so the updated text below reuses it, just
#     like the authored edit (sharedStrings count 28 -> 27).
$ws.Range("E35").ClearContents()
$ws.Range("A7").Value = "GM8724S009 (Lab1) (lab h bridge)"

# --- Selection moved to A8 ---
$ws.Range("A8").Select()

# --- Table 2 (rows 20-29): new stall-torque readings (column C); the
#     rad/s column D recalculates automatically from =C*0.10472 ---
$ws.Range("C21").Value = 164.02860000000001
$ws.Range("C22").Value = 266
$ws.Range("C23").Value = 287.43610000000001
$ws.Range("C24").Value = 129.08949999999999
$ws.Range("C25").Value = 140.01339999999999
$ws.Range("C26").Value = 234.69880000000001
$ws.Range("C27").Value = 199.1831
$ws.Range("C28").Value = 149.9992
$ws.Range("C29").Value = 120

# --- Table 3 (rows 33-41): new stall-torque readings (column C) ---
$ws.Range("C33").Value = 192.69280000000001
$ws.Range("C34").Value = 249.09780000000001
$ws.Range("C35").Value = 330.4289
$ws.Range("C36").Value = 136.9496
$ws.Range("C37").Value = 147.9128
$ws.Range("C38").Value = 209.69900000000001
$ws.Range("C39").Value = 174.99680000000001
$ws.Range("C40").Value = 131.25
$ws.Range("C41").Value = 105

# --- New (empty, formatted-only) helper cells to the right of table 3,
#     matching the widened used range (A1:K41) and the new G31:J31
#     merged header cell from the authored edit ---
$ws.Range("A7").Copy() | Out-Null
$ws.Range("G31:J31").PasteSpecial(-4104) | Out-Null
$ws.Range("G31:J31").Merge()

$ws.Range("C18").Copy() | Out-Null
$ws.Range("I32").PasteSpecial(-4104) | Out-Null
$ws.Range("G35:K35").PasteSpecial(-4104) | Out-Null
$ws.Range("G38").PasteSpecial(-4104) | Out-Null
$ws.Range("I38").PasteSpecial(-4104) | Out-Null
$ws.Range("G39").PasteSpecial(-4104) | Out-Null
$ws.Range("I39").PasteSpecial(-4104) | Out-Null
$ws.Range("I40").PasteSpecial(-4104) | Out-Null
$ws.Range("I41").PasteSpecial(-4104) | Out-Null

$ws.Range("B14").Copy() | Out-Null
$ws.Range("H38").PasteSpecial(-4104) | Out-Null
$ws.Range("H39").PasteSpecial(-4104) | Out-Null

$ws.Range("B18").Copy() | Out-Null
$ws.Range("H40").PasteSpecial(-4104) | Out-Null
$ws.Range("H41").PasteSpecial(-4104) | Out-Null

$excel.CutCopyMode = 0

Write-Output "edit applied"
